$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.2041522491349481
$ws.Cells.Item(2, 3).Value2 = 0.5501730103806228
$ws.Cells.Item(2, 10).Value2 = 0.01384083044982699
$ws.Cells.Item(2, 16).Value2 = 0.1453287197231834
$ws.Cells.Item(2, 19).Value2 = 0.08650519031141868
$ws.Cells.Item(3, 2).Value2 = 0.0124223602484472
$ws.Cells.Item(3, 10).Value2 = 0.02484472049689441
$ws.Cells.Item(3, 16).Value2 = 0.7701863354037267
$ws.Cells.Item(3, 19).Value2 = 0.1925465838509317
$ws.Cells.Item(4, 10).Value2 = 0.06818181818181818
$ws.Cells.Item(4, 16).Value2 = 0.6363636363636364
$ws.Cells.Item(4, 19).Value2 = 0.2954545454545455
$ws.Cells.Item(6, 2).Value2 = 0.0371900826446281
$ws.Cells.Item(6, 4).Value2 = 0.02066115702479339
$ws.Cells.Item(6, 6).Value2 = 0.06611570247933884
$ws.Cells.Item(6, 10).Value2 = 0.3264462809917356
$ws.Cells.Item(6, 15).Value2 = 0.01652892561983471
$ws.Cells.Item(6, 17).Value2 = 0.1487603305785124
$ws.Cells.Item(6, 18).Value2 = 0.06198347107438017
$ws.Cells.Item(6, 19).Value2 = 0.3223140495867768
$ws.Cells.Item(7, 2).Value2 = 0.07655502392344497
$ws.Cells.Item(7, 4).Value2 = 0.01913875598086124
$ws.Cells.Item(7, 5).Value2 = 0.004784688995215311
$ws.Cells.Item(7, 6).Value2 = 0.06220095693779904
$ws.Cells.Item(7, 10).Value2 = 0.1483253588516746
$ws.Cells.Item(7, 15).Value2 = 0.02392344497607655
$ws.Cells.Item(7, 17).Value2 = 0.1196172248803828
$ws.Cells.Item(7, 18).Value2 = 0.1291866028708134
$ws.Cells.Item(7, 19).Value2 = 0.4162679425837321
$ws.Cells.Item(8, 2).Value2 = 0.07061503416856492
$ws.Cells.Item(8, 4).Value2 = 0.01594533029612756
$ws.Cells.Item(8, 5).Value2 = 0.002277904328018223
$ws.Cells.Item(8, 6).Value2 = 0.03644646924829157
$ws.Cells.Item(8, 10).Value2 = 0.1230068337129841
$ws.Cells.Item(8, 15).Value2 = 0.01138952164009112
$ws.Cells.Item(8, 17).Value2 = 0.1708428246013667
$ws.Cells.Item(8, 18).Value2 = 0.1298405466970387
$ws.Cells.Item(8, 19).Value2 = 0.4396355353075171
$ws.Cells.Item(9, 2).Value2 = 0.1069767441860465
$ws.Cells.Item(9, 4).Value2 = 0.0186046511627907
$ws.Cells.Item(9, 6).Value2 = 0.06046511627906977
$ws.Cells.Item(9, 10).Value2 = 0.1162790697674419
$ws.Cells.Item(9, 15).Value2 = 0.02325581395348837
$ws.Cells.Item(9, 17).Value2 = 0.1534883720930233
$ws.Cells.Item(9, 18).Value2 = 0.1023255813953488
$ws.Cells.Item(9, 19).Value2 = 0.4186046511627907
$ws.Cells.Item(10, 2).Value2 = 0.1064446053584359
$ws.Cells.Item(10, 4).Value2 = 0.01737871107892831
$ws.Cells.Item(10, 6).Value2 = 0.07385952208544533
$ws.Cells.Item(10, 10).Value2 = 0.1194786386676322
$ws.Cells.Item(10, 15).Value2 = 0.02172338884866039
$ws.Cells.Item(10, 17).Value2 = 0.1976828385228095
$ws.Cells.Item(10, 18).Value2 = 0.0890658942795076
$ws.Cells.Item(10, 19).Value2 = 0.3743664011585807
$ws.Cells.Item(11, 7).Value2 = 0.1335403726708075
$ws.Cells.Item(11, 10).Value2 = 0.08074534161490683
$ws.Cells.Item(11, 11).Value2 = 0.2018633540372671
$ws.Cells.Item(11, 12).Value2 = 0.5652173913043478
$ws.Cells.Item(11, 19).Value2 = 0.01863354037267081
$ws.Cells.Item(12, 7).Value2 = 0.7564766839378239
$ws.Cells.Item(12, 10).Value2 = 0.1761658031088083
$ws.Cells.Item(12, 11).Value2 = 0.01036269430051814
$ws.Cells.Item(12, 12).Value2 = 0.03626943005181347
$ws.Cells.Item(12, 19).Value2 = 0.02072538860103627
$ws.Cells.Item(13, 7).Value2 = 0.6842105263157895
$ws.Cells.Item(13, 10).Value2 = 0.2368421052631579
$ws.Cells.Item(13, 19).Value2 = 0.07894736842105263
$ws.Cells.Item(15, 6).Value2 = 0.0411522633744856
$ws.Cells.Item(15, 8).Value2 = 0.1440329218106996
$ws.Cells.Item(15, 9).Value2 = 0.06584362139917696
$ws.Cells.Item(15, 10).Value2 = 0.3703703703703703
$ws.Cells.Item(15, 11).Value2 = 0.04938271604938271
$ws.Cells.Item(15, 13).Value2 = 0.00411522633744856
$ws.Cells.Item(15, 14).Value2 = 0.00411522633744856
$ws.Cells.Item(15, 15).Value2 = 0.06995884773662552
$ws.Cells.Item(15, 19).Value2 = 0.2510288065843622
$ws.Cells.Item(16, 6).Value2 = 0.02162162162162162
$ws.Cells.Item(16, 8).Value2 = 0.1891891891891892
$ws.Cells.Item(16, 9).Value2 = 0.08648648648648649
$ws.Cells.Item(16, 10).Value2 = 0.3837837837837838
$ws.Cells.Item(16, 11).Value2 = 0.1135135135135135
$ws.Cells.Item(16, 13).Value2 = 0.01081081081081081
$ws.Cells.Item(16, 15).Value2 = 0.03243243243243243
$ws.Cells.Item(16, 19).Value2 = 0.1621621621621622
$ws.Cells.Item(17, 6).Value2 = 0.02262443438914027
$ws.Cells.Item(17, 8).Value2 = 0.1719457013574661
$ws.Cells.Item(17, 9).Value2 = 0.1176470588235294
$ws.Cells.Item(17, 10).Value2 = 0.4072398190045249
$ws.Cells.Item(17, 11).Value2 = 0.09049773755656108
$ws.Cells.Item(17, 13).Value2 = 0.01809954751131222
$ws.Cells.Item(17, 15).Value2 = 0.05429864253393665
$ws.Cells.Item(17, 19).Value2 = 0.1176470588235294
$ws.Cells.Item(18, 6).Value2 = 0.03688524590163934
$ws.Cells.Item(18, 8).Value2 = 0.1598360655737705
$ws.Cells.Item(18, 9).Value2 = 0.1024590163934426
$ws.Cells.Item(18, 10).Value2 = 0.3647540983606558
$ws.Cells.Item(18, 11).Value2 = 0.09426229508196721
$ws.Cells.Item(18, 13).Value2 = 0.02049180327868852
$ws.Cells.Item(18, 15).Value2 = 0.07377049180327869
$ws.Cells.Item(18, 19).Value2 = 0.1475409836065574
$ws.Cells.Item(19, 6).Value2 = 0.01605839416058394
$ws.Cells.Item(19, 8).Value2 = 0.1912408759124088
$ws.Cells.Item(19, 9).Value2 = 0.0781021897810219
$ws.Cells.Item(19, 10).Value2 = 0.3897810218978102
$ws.Cells.Item(19, 11).Value2 = 0.1153284671532847
$ws.Cells.Item(19, 13).Value2 = 0.01897810218978102
$ws.Cells.Item(19, 15).Value2 = 0.072992700729927
$ws.Cells.Item(19, 19).Value2 = 0.1175182481751825
